$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employee ID
$ws.Range("C1").Value = "PCS0015"

# Update employee name
$ws.Range("C2").Value = "Mr. Santhosh Saravanan"

# Update monthly gross amounts
$ws.Range("D4").Value = 79200
$ws.Range("E4").Value = 99000
